$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bump the "Förändrad" (changed) date in column C from 2023-09-08 (45177)
# to 2023-09-09 (45178) for every data row (rows 2 through 118).
for ($r = 2; $r -le 118; $r++) {
    $ws.Cells.Item($r, 3).Value = 45178
}

# Row 118 gained an explicit row height in the new file.
$ws.Rows(118).RowHeight = 15

# Append the new row 119 with the new logging/felling notice.
$ws.Range("A119").Value = "A 41952-2023"

$ws.Range("B119").Value = 45177
$ws.Range("B119").NumberFormat = "YYYY-MM-DD"

$ws.Range("C119").Value = 45178
$ws.Range("C119").NumberFormat = "YYYY-MM-DD"

$ws.Range("D119").Value = "STOCKHOLMS LÄN"
$ws.Range("E119").Value = "SIGTUNA"
$ws.Range("F119").Value = "Allmännings- och besparingsskogar"
$ws.Range("G119").Value = 8.5
$ws.Range("H119:Q119").Value = 0

$ws.Range("R119").WrapText = $true
